# Apply targeted cell updates per the commit diff (regenerated site data).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1610
$ws.Range("F5").Value = 9193
$ws.Range("G6").Value = "不可售"
$ws.Range("G7").Value = 70
$ws.Range("F8").Value = 1280
$ws.Range("F11").Value = 603
$ws.Range("F13").Value = 158
$ws.Range("F17").Value = 1535
$ws.Range("F18").Value = 1335
$ws.Range("F21").Value = 1404
$ws.Range("F22").Value = 92
$ws.Range("F25").Value = 103
$ws.Range("F28").Value = 328
$ws.Range("F29").Value = 328
$ws.Range("F32").Value = 34
$ws.Range("F34").Value = 216
$ws.Range("F36").Value = 583
$ws.Range("F37").Value = 614
$ws.Range("F39").Value = 140
$ws.Range("F41").Value = 159
$ws.Range("F42").Value = 96
$ws.Range("F43").Value = 519
$ws.Range("F44").Value = 1239
$ws.Range("F45").Value = 696
$ws.Range("F46").Value = 237
$ws.Range("F48").Value = 47

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 58
$ws.Range("F14").Value = 2
$ws.Range("F16").Value = 674
$ws.Range("F21").Value = 8
$ws.Range("F24").Value = 939
$ws.Range("F25").Value = 23
$ws.Range("F26").Value = 1038
$ws.Range("F27").Value = 247
$ws.Range("F28").Value = 641
$ws.Range("F30").Value = 246

# --- Sheet 3: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 325
$ws.Range("F6").Value = 148
$ws.Range("F7").Value = 2182
$ws.Range("F8").Value = 3263
$ws.Range("F9").Value = 31

# --- Sheet 4: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1610
$ws.Range("F6").Value = 9193
$ws.Range("F7").Value = 325
$ws.Range("F8").Value = 148
$ws.Range("C9").Value = "上海·EVANGELION× PrismLand · 新世纪福音战士官方授权主题店"
$ws.Range("D9").Value = "南京东路830号第一百货商业中心B馆5楼(海底捞旁边) 第一百货商业中心"
$ws.Range("E9").Value = "2024.05.25 00:00-07.22 23:59"
$ws.Range("F9").Value = 2182
$ws.Range("G9").Value = 20
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=85030"
$ws.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202404/K3kIpfaB1714445776157.jpeg"
$ws.Range("C10").Value = "上海·「排球少年!!垃圾场决战 × animate cafe」"
$ws.Range("D10").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws.Range("E10").Value = "2024.05.25 00:00-07.02 23:59"
$ws.Range("F10").Value = 3263
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85283"
$ws.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202405/vy2vecK11715162037223.jpeg"
$ws.Range("C11").Value = "上海·拉帮结派ONLY"
$ws.Range("D11").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E11").Value = "2024.05.25 14:00-05.25 19:00"
$ws.Range("F11").Value = 124
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=85091"
$ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202405/Gqv3tfiB1714795562310.jpeg"
$ws.Range("C12").Value = "上海·第六届Redamancy动漫游戏嘉年华"
$ws.Range("D12").Value = "中山北路3300号4楼 上海环球港"
$ws.Range("E12").Value = "2024.05.25 10:00-05.26 17:00"
$ws.Range("F12").Value = 1280
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=84632"
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202404/im8C39eo1713190504331.png"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2024-05-26"
$ws.Range("C13").Value = "上海·五十岚隼士&小池亮介·2024见面会"
$ws.Range("D13").Value = "长寿路街道万航渡后路19号 上海瓦肆文化传播有限公司"
$ws.Range("E13").Value = "2024.05.26 11:00-05.26 15:30"
$ws.Range("F13").Value = 162
$ws.Range("G13").Value = 480
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=84615"
$ws.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202404/LwpFUbIJ1713857706981.jpeg"
$ws.Range("F16").Value = 603
$ws.Range("F17").Value = 158
$ws.Range("F19").Value = 1535
$ws.Range("F20").Value = 674
$ws.Range("F21").Value = 1335
$ws.Range("F23").Value = 31
$ws.Range("F24").Value = 1404
$ws.Range("F26").Value = 103
$ws.Range("F28").Value = 328
$ws.Range("F29").Value = 328
$ws.Range("F33").Value = 34
$ws.Range("F34").Value = 939
$ws.Range("F36").Value = 216
$ws.Range("F37").Value = 1038
$ws.Range("F38").Value = 247
$ws.Range("F39").Value = 583
$ws.Range("F40").Value = 614
$ws.Range("F41").Value = 641
$ws.Range("F42").Value = 140
$ws.Range("F43").Value = 246
$ws.Range("F44").Value = 159
$ws.Range("F45").Value = 96
$ws.Range("F47").Value = 519
$ws.Range("F48").Value = 696
